$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)
$nl = [char]11

$tbl.Cell(1, 1).Range.Text = "96 x 52" + $nl + "  5    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "6|    |"
$tbl.Cell(1, 2).Range.Text = "69 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "6|    |" + $nl + "9|    |"
$tbl.Cell(1, 3).Range.Text = "16 x 87" + $nl + "  8    7" + $nl + "  ----" + $nl + "1|    |" + $nl + "6|    |"
$tbl.Cell(2, 1).Range.Text = "24 x 98" + $nl + "  9    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$tbl.Cell(2, 2).Range.Text = "59 x 82" + $nl + "  8    2" + $nl + "  ----" + $nl + "5|    |" + $nl + "9|    |"
$tbl.Cell(2, 3).Range.Text = "80 x 39" + $nl + "  3    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "0|    |"
$tbl.Cell(3, 1).Range.Text = "10 x 50" + $nl + "  5    0" + $nl + "  ----" + $nl + "1|    |" + $nl + "0|    |"
$tbl.Cell(3, 2).Range.Text = "43 x 15" + $nl + "  1    5" + $nl + "  ----" + $nl + "4|    |" + $nl + "3|    |"
$tbl.Cell(3, 3).Range.Text = "11 x 96" + $nl + "  9    6" + $nl + "  ----" + $nl + "1|    |" + $nl + "1|    |"
$tbl.Cell(4, 1).Range.Text = "26 x 14" + $nl + "  1    4" + $nl + "  ----" + $nl + "2|    |" + $nl + "6|    |"
$tbl.Cell(4, 2).Range.Text = "17 x 41" + $nl + "  4    1" + $nl + "  ----" + $nl + "1|    |" + $nl + "7|    |"
$tbl.Cell(4, 3).Range.Text = "46 x 40" + $nl + "  4    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "6|    |"
$tbl.Cell(5, 1).Range.Text = "21 x 58" + $nl + "  5    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "1|    |"
$tbl.Cell(5, 2).Range.Text = "60 x 72" + $nl + "  7    2" + $nl + "  ----" + $nl + "6|    |" + $nl + "0|    |"
$tbl.Cell(5, 3).Range.Text = "29 x 27" + $nl + "  2    7" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
